# The "Training Dashboard" report was refreshed: the "PERIOD TO EXPIRE"
# (H) and "LAST UPDATE" (I) columns move forward by one refresh cycle
# (08-Sep-2025 -> 16-Sep-2025, i.e. -8 days remaining), and the header
# row is restyled with a white, bold font so it reads clearly against
# its dark-blue fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Data refresh: PERIOD TO EXPIRE (H) and LAST UPDATE (I) ----------------
# New "last update" date for every training row.
$newUpdateDate = "16-Sep-2025"

# Keep the LAST UPDATE column as literal text (not an Excel date serial) by
# writing it as Text before assigning, matching how the report already
# stores its dates as plain strings.
$ws.Range("I3:I19").NumberFormat = "@"

$periods = @{
    3  = 590
    4  = 591
    5  = 594
    6  = 594
    7  = 590
    8  = 590
    9  = 591
    10 = 594
    11 = 594
    12 = 591
    13 = 594
    14 = -19618
    15 = 278
    16 = 313
    17 = 313
    18 = 313
    19 = 308
}

foreach ($row in $periods.Keys) {
    $ws.Cells.Item($row, 8).Value = $periods[$row]
    $ws.Cells.Item($row, 9).Value = $newUpdateDate
}

# --- Header row styling: bold white text on the dark-blue fill -------------
$ws.Range("A2:K2").Font.Color = 16777215
